# edit.ps1 -- "added harvard case classification"
#
# This sheet rolls up per-app evaluation stats (average / variance / std Dev)
# for several symptom-checker apps plus a panel of human-doctor raters
# (doctor_MA, doctor_MA_old, doctor_NJ, doctor_NJ_old, doctor_TH, doctor_TH_old),
# with BP/BQ holding the cross-doctor rollup ("average_doctor").
#
# Adding the Harvard case classification triggers a recompute of those stats:
#   - BP1/BQ1 swap meaning: BP ("average_doctor") becomes "average_doctor_old",
#     and BQ ("average_doctor_old") becomes the new "average_doctor".
#   - The data follows the label: each row's old BP value (previously
#     "average_doctor") is carried over into BQ, and BP is filled with the
#     freshly recomputed average.
#   - The other *_old columns (Ada_old, Avey_old, Babylon_old, Buoy_old,
#     K health_old, WebMD_old) and the doctor_* rater columns are updated in
#     place with the newly recomputed figures, for data rows 4-13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 1: relabel the average-doctor columns ---
$ws.Range("BP1").Value = "average_doctor_old"
$ws.Range("BQ1").Value = "average_doctor"

# --- Recomputed statistic values, rows 4-13 ---
# Row 4
$ws.Range("E4").Value = 0.436
$ws.Range("F4").Value = 0.07000000000000001
$ws.Range("G4").Value = 0.265
$ws.Range("N4").Value = 0.439
$ws.Range("O4").Value = 0.065
$ws.Range("P4").Value = 0.255
$ws.Range("Q4").Value = 0.019
$ws.Range("R4").Value = 0.014
$ws.Range("S4").Value = 0.118
$ws.Range("W4").Value = 0.297
$ws.Range("X4").Value = 0.114
$ws.Range("Y4").Value = 0.337
$ws.Range("AI4").Value = 0.335
$ws.Range("AJ4").Value = 0.08699999999999999
$ws.Range("AK4").Value = 0.294
$ws.Range("AU4").Value = 0.194
$ws.Range("AV4").Value = 0.03
$ws.Range("AW4").Value = 0.174
$ws.Range("BA4").Value = 2.016
$ws.Range("BB4").Value = 0.152
$ws.Range("BC4").Value = 0.39
$ws.Range("BG4").Value = 0.733
$ws.Range("BH4").Value = 0.135
$ws.Range("BI4").Value = 0.367
$ws.Range("BM4").Value = 0.723
$ws.Range("BN4").Value = 0.076
$ws.Range("BO4").Value = 0.276
$ws.Range("BP4").Value = 0.672
$ws.Range("BQ4").Value = 0.707
# Row 5
$ws.Range("E5").Value = 0.55
$ws.Range("F5").Value = 0.075
$ws.Range("G5").Value = 0.273
$ws.Range("N5").Value = 0.734
$ws.Range("O5").Value = 0.079
$ws.Range("P5").Value = 0.28
$ws.Range("Q5").Value = 0.01
$ws.Range("R5").Value = 0.003
$ws.Range("S5").Value = 0.05
$ws.Range("W5").Value = 0.278
$ws.Range("X5").Value = 0.105
$ws.Range("Y5").Value = 0.324
$ws.Range("AI5").Value = 0.359
$ws.Range("AJ5").Value = 0.097
$ws.Range("AK5").Value = 0.311
$ws.Range("AU5").Value = 0.373
$ws.Range("AV5").Value = 0.099
$ws.Range("AW5").Value = 0.315
$ws.Range("BA5").Value = 1.335
$ws.Range("BB5").Value = 0.077
$ws.Range("BC5").Value = 0.278
$ws.Range("BG5").Value = 0.395
$ws.Range("BH5").Value = 0.046
$ws.Range("BI5").Value = 0.214
$ws.Range("BM5").Value = 0.543
$ws.Range("BN5").Value = 0.057
$ws.Range("BO5").Value = 0.239
$ws.Range("BP5").Value = 0.445
$ws.Range("BQ5").Value = 0.456
# Row 6
$ws.Range("E6").Value = 0.486
$ws.Range("N6").Value = 0.549
$ws.Range("Q6").Value = 0.013
$ws.Range("W6").Value = 0.287
$ws.Range("AI6").Value = 0.347
$ws.Range("AU6").Value = 0.255
$ws.Range("BA6").Value = 1.598
$ws.Range("BG6").Value = 0.513
$ws.Range("BM6").Value = 0.62
$ws.Range("BP6").Value = 0.533
$ws.Range("BQ6").Value = 0.551
# Row 7
$ws.Range("E7").Value = 0.523
$ws.Range("N7").Value = 0.647
$ws.Range("Q7").Value = 0.011
$ws.Range("W7").Value = 0.282
$ws.Range("AI7").Value = 0.354
$ws.Range("AU7").Value = 0.315
$ws.Range("BA7").Value = 1.428
$ws.Range("BG7").Value = 0.435
$ws.Range("BM7").Value = 0.571
$ws.Range("BP7").Value = 0.476
$ws.Range("BQ7").Value = 0.489
# Row 8
$ws.Range("E8").Value = 0.617
$ws.Range("F8").Value = 0.101
$ws.Range("G8").Value = 0.318
$ws.Range("N8").Value = 0.771
$ws.Range("O8").Value = 0.066
$ws.Range("P8").Value = 0.257
$ws.Range("Q8").Value = 0.01
$ws.Range("S8").Value = 0.076
$ws.Range("W8").Value = 0.318
$ws.Range("X8").Value = 0.123
$ws.Range("Y8").Value = 0.351
$ws.Range("AI8").Value = 0.383
$ws.Range("AJ8").Value = 0.129
$ws.Range("AK8").Value = 0.359
$ws.Range("AU8").Value = 0.321
$ws.Range("AV8").Value = 0.08699999999999999
$ws.Range("AW8").Value = 0.294
$ws.Range("BA8").Value = 1.737
$ws.Range("BB8").Value = 0.121
$ws.Range("BC8").Value = 0.347
$ws.Range("BG8").Value = 0.5580000000000001
$ws.Range("BH8").Value = 0.103
$ws.Range("BI8").Value = 0.32
$ws.Range("BM8").Value = 0.6899999999999999
$ws.Range("BN8").Value = 0.062
$ws.Range("BO8").Value = 0.249
$ws.Range("BP8").Value = 0.579
$ws.Range("BQ8").Value = 0.599
# Row 9
$ws.Range("E9").Value = 0.556
$ws.Range("F9").Value = 0.247
$ws.Range("G9").Value = 0.497
$ws.Range("N9").Value = 0.667
$ws.Range("O9").Value = 0.222
$ws.Range("P9").Value = 0.471
$ws.Range("W9").Value = 0.222
$ws.Range("X9").Value = 0.173
$ws.Range("Y9").Value = 0.416
$ws.Range("AI9").Value = 0.296
$ws.Range("AJ9").Value = 0.209
$ws.Range("AK9").Value = 0.457
$ws.Range("BA9").Value = 1.666
$ws.Range("BB9").Value = 0.247
$ws.Range("BC9").Value = 0.497
$ws.Range("BG9").Value = 0.58
$ws.Range("BH9").Value = 0.244
$ws.Range("BI9").Value = 0.494
$ws.Range("BM9").Value = 0.642
$ws.Range("BN9").Value = 0.23
$ws.Range("BO9").Value = 0.479
$ws.Range("BP9").Value = 0.555
$ws.Range("BQ9").Value = 0.5659999999999999
# Row 10
$ws.Range("E10").Value = 0.6909999999999999
$ws.Range("F10").Value = 0.213
$ws.Range("G10").Value = 0.462
$ws.Range("N10").Value = 0.877
$ws.Range("O10").Value = 0.108
$ws.Range("P10").Value = 0.329
$ws.Range("W10").Value = 0.395
$ws.Range("X10").Value = 0.239
$ws.Range("Y10").Value = 0.489
$ws.Range("AI10").Value = 0.42
$ws.Range("AJ10").Value = 0.244
$ws.Range("AK10").Value = 0.494
$ws.Range("AU10").Value = 0.321
$ws.Range("AV10").Value = 0.218
$ws.Range("AW10").Value = 0.467
$ws.Range("BA10").Value = 2.074
$ws.Range("BB10").Value = 0.244
$ws.Range("BC10").Value = 0.494
$ws.Range("BG10").Value = 0.642
$ws.Range("BH10").Value = 0.23
$ws.Range("BI10").Value = 0.479
$ws.Range("BM10").Value = 0.852
$ws.Range("BN10").Value = 0.126
$ws.Range("BO10").Value = 0.355
$ws.Range("BP10").Value = 0.6909999999999999
$ws.Range("BQ10").Value = 0.714
# Row 11
$ws.Range("E11").Value = 0.728
$ws.Range("F11").Value = 0.198
$ws.Range("G11").Value = 0.445
$ws.Range("N11").Value = 0.889
$ws.Range("O11").Value = 0.099
$ws.Range("P11").Value = 0.314
$ws.Range("W11").Value = 0.395
$ws.Range("X11").Value = 0.239
$ws.Range("Y11").Value = 0.489
$ws.Range("AI11").Value = 0.457
$ws.Range("AJ11").Value = 0.248
$ws.Range("AK11").Value = 0.498
$ws.Range("AU11").Value = 0.444
$ws.Range("AV11").Value = 0.247
$ws.Range("AW11").Value = 0.497
$ws.Range("BA11").Value = 2.074
$ws.Range("BB11").Value = 0.244
$ws.Range("BC11").Value = 0.494
$ws.Range("BG11").Value = 0.642
$ws.Range("BH11").Value = 0.23
$ws.Range("BI11").Value = 0.479
$ws.Range("BM11").Value = 0.852
$ws.Range("BN11").Value = 0.126
$ws.Range("BO11").Value = 0.355
$ws.Range("BP11").Value = 0.6909999999999999
$ws.Range("BQ11").Value = 0.717
# Row 12
$ws.Range("E12").Value = 1.441
$ws.Range("F12").Value = 0.823
$ws.Range("G12").Value = 0.907
$ws.Range("N12").Value = 1.473
$ws.Range("O12").Value = 1.006
$ws.Range("P12").Value = 1.003
$ws.Range("W12").Value = 1.594
$ws.Range("X12").Value = 0.554
$ws.Range("Y12").Value = 0.744
$ws.Range("AI12").Value = 1.703
$ws.Range("AJ12").Value = 1.29
$ws.Range("AK12").Value = 1.136
$ws.Range("AU12").Value = 2.711
$ws.Range("AV12").Value = 2.785
$ws.Range("AW12").Value = 1.669
$ws.Range("BA12").Value = 3.774
$ws.Range("BB12").Value = 0.437
$ws.Range("BC12").Value = 0.661
$ws.Range("BG12").Value = 1.115
$ws.Range("BH12").Value = 0.141
$ws.Range("BI12").Value = 0.375
$ws.Range("BM12").Value = 1.319
$ws.Range("BN12").Value = 0.362
$ws.Range("BO12").Value = 0.602
$ws.Range("BP12").Value = 1.258
$ws.Range("BQ12").Value = 1.28
# Row 13
$ws.Range("E13").Value = 1.554
$ws.Range("F13").Value = 0.615
$ws.Range("G13").Value = 0.784
$ws.Range("N13").Value = 2.007
$ws.Range("O13").Value = 0.781
$ws.Range("P13").Value = 0.884
$ws.Range("W13").Value = 1.024
$ws.Range("X13").Value = 0.192
$ws.Range("Y13").Value = 0.438
$ws.Range("AI13").Value = 1.268
$ws.Range("AJ13").Value = 0.372
$ws.Range("AK13").Value = 0.61
$ws.Range("AU13").Value = 2.22
$ws.Range("AV13").Value = 0.593
$ws.Range("AW13").Value = 0.77
$ws.Range("BA13").Value = 2.318
$ws.Range("BB13").Value = 0.294
$ws.Range("BC13").Value = 0.543
$ws.Range("BG13").Value = 0.582
$ws.Range("BH13").Value = 0.07199999999999999
$ws.Range("BI13").Value = 0.268
$ws.Range("BM13").Value = 0.873
$ws.Range("BN13").Value = 0.271
$ws.Range("BO13").Value = 0.521
$ws.Range("BP13").Value = 0.773
$ws.Range("BQ13").Value = 0.719
